# Denmark 3rd-division 2023-2024: reorder a handful of match rows (betting
# odds columns F:V swapped between paired rows) and append the newly
# scraped Avarta vs Vejgaard match as row 89.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap columns F:V between the following row pairs ------------------
$pairs = @(
    @(64, 65),
    @(68, 69),
    @(70, 71),
    @(76, 77),
    @(80, 82)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("F$($r1):V$($r1)")
    $range2 = $ws.Range("F$($r2):V$($r2)")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}

# --- Append new row 89 (Avarta vs Vejgaard, 12/11/2023) -----------------
# Copy formatting (number formats / styles) from the last existing row.
$ws.Range("A88:V88").Copy()
$ws.Range("A89:V89").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A89").Value = 88
$ws.Range("B89").Value = "denmark"
$ws.Range("C89").Value = "3rd-division"
$ws.Range("D89").Value = "2023-2024"
$ws.Range("E89").Value = 45242.54166666666
$ws.Range("F89").Value = "Avarta"
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = "Vejgaard"
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2.39
$ws.Range("K89").Value = "11/11/2023 01:12"
$ws.Range("L89").Value = 2.62
$ws.Range("M89").Value = "12/11/2023 12:40"
$ws.Range("N89").Value = 3.2
$ws.Range("O89").Value = "11/11/2023 01:12"
$ws.Range("P89").Value = 3.21
$ws.Range("Q89").Value = "12/11/2023 12:46"
$ws.Range("R89").Value = 2.53
$ws.Range("S89").Value = "11/11/2023 01:12"
$ws.Range("T89").Value = 2.55
$ws.Range("U89").Value = "12/11/2023 12:40"
$ws.Range("V89").Value = "https://www.betexplorer.com/football/denmark/3rd-division/avarta-vejgaard/Gp8PaUb7/"
